# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@8e4a450c507ef6f746e072652acbb72e9504f19a
# Update the "Metadata" sheet (version bump, new publish date, publisher, jurisdiction)
# and the "Elements" sheet (regenerated Short/Definition text for the root element),
# reflecting the StructureDefinition-canonical-measure being regenerated against the
# new IG build.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# The old sheet had a duplicated "Contact" row (rows 10 & 11) caused by a renderer bug
# for ContactDetail. The regenerated IG collapses that into a single "Jurisdiction" row.
$meta.Rows.Item(11).Delete()

# Version bump
$meta.Range("B3").Value = "6.0.0"

# New publication date/time
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a display value
$meta.Range("B9").Value = "Alvearie Team"

# Former duplicate "Contact" row now shows Jurisdiction info
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Elements sheet: root Extension row's Short/Definition now reflect the
# canonical-measure extension specifically instead of the generic placeholder text.
$elements.Range("K2").Value = "CanonicalMeasure"
$elements.Range("L2").Value = "Canonical reference to the specific version of the measure used to generate the resource."
